# Insert a new weekly price record at row 86 ("Fruta / hortaliza, semanal").
# This shifts the existing rows 86..191 down to 87..192 and fills the newly
# opened row 86 with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(86).Insert()

$ws.Cells.Item(86, 1).Value = 1
$ws.Cells.Item(86, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(86, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(86, 4).Value = 44467
$ws.Cells.Item(86, 5).Value = 15
$ws.Cells.Item(86, 6).Value = 100114013
$ws.Cells.Item(86, 7).Value = "Zanahoria"
$ws.Cells.Item(86, 8).Value = "Sin especificar"
$ws.Cells.Item(86, 9).Value = "Primera"
$ws.Cells.Item(86, 10).Value = 100
$ws.Cells.Item(86, 11).Value = 6500
$ws.Cells.Item(86, 12).Value = 7000
$ws.Cells.Item(86, 13).Value = 6750
$ws.Cells.Item(86, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(86, 15).Value = "Valle de Camiña"
$ws.Cells.Item(86, 16).Value = 270
$ws.Cells.Item(86, 17).Value = 25
$ws.Cells.Item(86, 18).Value = "Hortaliza"
